$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text first so numeric-looking strings (e.g. "577.57")
# are stored as literal text instead of being parsed as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "63.168.69"
$ws.Range("E2").Value = "  +0.21%  "
$ws.Range("D3").Value = "2.476.53"
$ws.Range("E3").Value = "  +2.42%  "
$ws.Range("E4").Value = "  -0.40%  "
$ws.Range("D5").Value = "577.57"
$ws.Range("E5").Value = "  +0.61%  "
$ws.Range("D6").Value = "146.91"
$ws.Range("E6").Value = "  +0.47%  "
$ws.Range("E8").Value = "  -0.35%  "
$ws.Range("D9").Value = "2.476.75"
$ws.Range("E9").Value = "  +1.06%  "
$ws.Range("E10").Value = "  +0.47%  "
$ws.Range("E11").Value = "  +1.61%  "
$ws.Range("D12").Value = "5.28"
$ws.Range("E12").Value = "  +0.61%  "
$ws.Range("E13").Value = "  +0.26%  "
$ws.Range("D14").Value = "28.69"
$ws.Range("E14").Value = "  +4.87%  "
$ws.Range("E15").Value = "  +0.75%  "
$ws.Range("D16").Value = "2.927.61"
$ws.Range("E16").Value = "  +2.43%  "
$ws.Range("D17").Value = "63.170.42"
$ws.Range("E17").Value = "  +0.61%  "
$ws.Range("D18").Value = "2.479.14"
$ws.Range("E18").Value = "  +1.45%  "
$ws.Range("D19").Value = "8.27"
$ws.Range("E19").Value = "  +4.24%  "
$ws.Range("E20").Value = "  +0.86%  "
$ws.Range("D21").Value = "329.65"
$ws.Range("E21").Value = "  +0.22%  "
$ws.Range("E22").Value = "  +10.24%  "
$ws.Range("D23").Value = "4.14"
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").Value = "  +0.22%  "
$ws.Range("D25").Value = "66.36"
$ws.Range("E25").Value = "  +1.12%  "
$ws.Range("D26").Value = "672.86"
$ws.Range("E26").Value = "  +5.82%  "
$ws.Range("D27").Value = "9.69"
$ws.Range("E27").Value = "  +13.61%  "
$ws.Range("E28").Value = "  +0.56%  "
$ws.Range("D29").Value = "2.619.39"
$ws.Range("E29").Value = "  +3.20%  "
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  -9.12%  "
$ws.Range("E31").Value = "  +2.68%  "
$ws.Range("D32").Value = "8.06"
$ws.Range("E32").Value = "  -2.19%  "
$ws.Range("E33").Value = "  +1.12%  "
$ws.Range("E34").Value = "  -3.45%  "
$ws.Range("E35").Value = "  +3.83%  "
$ws.Range("E36").Value = "  +0.34%  "
$ws.Range("D37").Value = "4.80"
$ws.Range("E37").Value = "  +0.76%  "
$ws.Range("D38").Value = "5.50"
$ws.Range("E38").Value = "  +1.28%  "
$ws.Range("E39").Value = "  -0.76%  "
$ws.Range("E40").Value = "  +0.67%  "
$ws.Range("D41").Value = "151.02"
$ws.Range("E41").Value = "  -1.31%  "
$ws.Range("E42").Value = "  -1.83%  "
$ws.Range("E43").Value = "  -0.40%  "
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("D45").Value = "0.0₆0313"
$ws.Range("E45").Value = "  +9.24%  "
$ws.Range("D46").Value = "154.45"
$ws.Range("E46").Value = "  +6.33%  "
$ws.Range("E47").Value = "  +15.94%  "
$ws.Range("D48").Value = "3.61"
$ws.Range("E48").Value = "  +0.06%  "
$ws.Range("D49").Value = "20.67"
$ws.Range("E49").Value = "  +0.66%  "
$ws.Range("E50").Value = "  +0.70%  "
$ws.Range("D51").Value = "0.0514"

# Restore the default (Normal) style on column D so we do not leave behind
# an extra text-format style that was not present in the original workbook.
$ws.Range("D2:D51").Style = "Normal"
